# FAST_holdings.xlsx: fill in the "Percent Change" column (E) with the
# day's computed percent-change values. The sheet ships protected, so it
# must be unprotected before the cells can be written, then re-protected
# to restore the original "this sheet is protected" state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

$ws.Range("E2").Value  = -0.02183502426113826
$ws.Range("E3").Value  = -0.003952109729164266
$ws.Range("E4").Value  = 0.004980301791421837
$ws.Range("E5").Value  = 0.0009242144177448175
$ws.Range("E6").Value  = 0.006039457790900427
$ws.Range("E8").Value  = 0.002238388360380483
$ws.Range("E9").Value  = -0.0002031006702322147
$ws.Range("E10").Value = -0.001791061378998982

$ws.Protect()
